# "further cleaning to metadata"
# - libraryProtocol (col K) is unified to a single value "E7420" (renamed
#   from "E7760") for every data row, replacing the old row-unique
#   E7761..E7766 placeholders used on rows 22-27.
# - roboticLibraryPrep (col L) switches from a bare boolean literal to an
#   explicit "=FALSE()" formula, with its own (new) font/number format.
# - The active selection moves from the L column over to the K column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 27

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $kCell = $ws.Cells.Item($r, 11)   # column K - libraryProtocol
    $lCell = $ws.Cells.Item($r, 12)   # column L - roboticLibraryPrep

    # Unify every libraryProtocol cell on the renamed value, and give it
    # the boolean-ish number format that column L used to carry.
    $kCell.NumberFormat = '"TRUE";"TRUE";"FALSE"'
    $kCell.Value = "E7420"

    # roboticLibraryPrep becomes a real formula instead of a static boolean.
    # Number format is reset to General *before* touching the font so no
    # transient (new-font + old-boolean-format) style gets interned.
    $lCell.NumberFormat = "General"
    $lCell.Formula = "=FALSE()"
    $lCell.Font.Name = "Arial"
    $lCell.Font.Size = 11
    $lCell.Font.Color = 0
}

$ws.Range("K2:K27").Select()
